$wb = $excel.ActiveWorkbook

# This script applies updated profit-calculation values to the Midgardsormr
# Profits workbook, sheet by sheet, as produced by the scheduled data refresh.

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 167.33333
$ws.Range("I11").Value = 167.33333
$ws.Range("K11").Value = 167.33333
$ws.Range("M11").Value = -27.33332999999999
$ws.Range("H17").Value = 1894.8823
$ws.Range("J17").Value = 1991.6
$ws.Range("L17").Value = 5974.799999999999
$ws.Range("N17").Value = -6310.799999999999
$ws.Range("H32").Value = 13426.308
$ws.Range("I32").Value = 19799.25
$ws.Range("K32").Value = 19799.25
$ws.Range("M32").Value = -19473.25
$ws.Range("H100").Value = 54436.75
$ws.Range("I100").Value = 76989.875
$ws.Range("K100").Value = 76989.875
$ws.Range("M100").Value = -76448.875
$ws.Range("H103").Value = 2115.0667
$ws.Range("I103").Value = 748.5
$ws.Range("J103").Value = 3026.111
$ws.Range("K103").Value = 2245.5
$ws.Range("L103").Value = 9078.332999999999
$ws.Range("M103").Value = -1659.5
$ws.Range("N103").Value = -10250.333
$ws.Range("H141").Value = 1408.8379
$ws.Range("I141").Value = 1066.0303
$ws.Range("K141").Value = 3198.0909
$ws.Range("M141").Value = 1981.9091

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 2592.75
$ws.Range("I31").Value = 2592.75
$ws.Range("K31").Value = 2592.75
$ws.Range("M31").Value = -2298.75
$ws.Range("H32").Value = 14413.083
$ws.Range("I32").Value = 14786.42
$ws.Range("K32").Value = 14786.42
$ws.Range("M32").Value = -14499.42
$ws.Range("H45").Value = 4354.3076
$ws.Range("I45").Value = 2042
$ws.Range("J45").Value = 9557
$ws.Range("K45").Value = 2042
$ws.Range("L45").Value = 9557
$ws.Range("M45").Value = -1665
$ws.Range("N45").Value = -10311
$ws.Range("H61").Value = 2933.17
$ws.Range("I61").Value = 776.8158
$ws.Range("K61").Value = 776.8158
$ws.Range("M61").Value = -564.8158
$ws.Range("H63").Value = 3409.1052
$ws.Range("I63").Value = 2027.3
$ws.Range("J63").Value = 4944.4443
$ws.Range("K63").Value = 2027.3
$ws.Range("L63").Value = 4944.4443
$ws.Range("M63").Value = -1341.3
$ws.Range("N63").Value = -6316.4443
$ws.Range("H66").Value = 3409.1052
$ws.Range("I66").Value = 2027.3
$ws.Range("J66").Value = 4944.4443
$ws.Range("K66").Value = 10136.5
$ws.Range("L66").Value = 24722.2215
$ws.Range("M66").Value = -6704.5
$ws.Range("N66").Value = -31586.2215
$ws.Range("H74").Value = 161695.53
$ws.Range("I74").Value = 207994.38
$ws.Range("K74").Value = 207994.38
$ws.Range("M74").Value = -207120.38
$ws.Range("H77").Value = 161695.53
$ws.Range("I77").Value = 207994.38
$ws.Range("K77").Value = 1039971.9
$ws.Range("M77").Value = -1035603.9
$ws.Range("H102").Value = 7827.282
$ws.Range("I102").Value = 8070.4062
$ws.Range("K102").Value = 8070.4062
$ws.Range("M102").Value = -6448.4062
$ws.Range("H132").Value = 1881.6052
$ws.Range("I132").Value = 1709.6
$ws.Range("J132").Value = 2072.7222
$ws.Range("K132").Value = 5128.799999999999
$ws.Range("L132").Value = 6218.1666
$ws.Range("M132").Value = -2598.799999999999
$ws.Range("N132").Value = -11278.1666
$ws.Range("H134").Value = 84857
$ws.Range("J134").Value = 84857
$ws.Range("L134").Value = 84857
$ws.Range("N134").Value = -94997
$ws.Range("H136").Value = 2933.17
$ws.Range("I136").Value = 776.8158
$ws.Range("K136").Value = 2330.4474
$ws.Range("M136").Value = 219.5526

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 19958.578
$ws.Range("I20").Value = 23985.379
$ws.Range("J20").Value = 1334.625
$ws.Range("K20").Value = 23985.379
$ws.Range("L20").Value = 1334.625
$ws.Range("M20").Value = -23738.379
$ws.Range("N20").Value = -1828.625
$ws.Range("H99").Value = 2852.9722
$ws.Range("I99").Value = 2675.577
$ws.Range("K99").Value = 2675.577
$ws.Range("M99").Value = -1177.577
$ws.Range("H134").Value = 3786.8157
$ws.Range("I134").Value = 1956.68
$ws.Range("K134").Value = 5870.04
$ws.Range("M134").Value = -3335.04

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3847805.8
$ws.Range("I31").Value = 4167935
$ws.Range("J31").Value = 6257
$ws.Range("K31").Value = 4167935
$ws.Range("L31").Value = 6257
$ws.Range("M31").Value = -4167640
$ws.Range("N31").Value = -6847
$ws.Range("H34").Value = 3847805.8
$ws.Range("I34").Value = 4167935
$ws.Range("J34").Value = 6257
$ws.Range("K34").Value = 4167935
$ws.Range("L34").Value = 6257
$ws.Range("M34").Value = -4167733
$ws.Range("N34").Value = -6661
$ws.Range("H44").Value = 1000
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("H69").Value = 26500
$ws.Range("I69").Value = 6750
$ws.Range("J69").Value = 66000
$ws.Range("K69").Value = 6750
$ws.Range("L69").Value = 66000
$ws.Range("M69").Value = -6001
$ws.Range("N69").Value = -67498
$ws.Range("H72").Value = 26500
$ws.Range("I72").Value = 6750
$ws.Range("J72").Value = 66000
$ws.Range("K72").Value = 20250
$ws.Range("L72").Value = 198000
$ws.Range("M72").Value = -16506
$ws.Range("N72").Value = -205488
$ws.Range("H107").Value = 1097.2307
$ws.Range("I107").Value = 381
$ws.Range("J107").Value = 1932.8334
$ws.Range("K107").Value = 381
$ws.Range("L107").Value = 1932.8334
$ws.Range("M107").Value = 1539
$ws.Range("N107").Value = -5772.8334
$ws.Range("H108").Value = 45155.2
$ws.Range("J108").Value = 48592
$ws.Range("L108").Value = 48592
$ws.Range("N108").Value = -56272
$ws.Range("H134").Value = 1573
$ws.Range("I134").Value = 1481.1428
$ws.Range("J134").Value = 1799.9412
$ws.Range("K134").Value = 4443.428400000001
$ws.Range("L134").Value = 5399.8236
$ws.Range("M134").Value = -1908.428400000001
$ws.Range("N134").Value = -10469.8236
$ws.Range("M44").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7223.4375
$ws.Range("I56").Value = 7223.4375
$ws.Range("K56").Value = 7223.4375
$ws.Range("M56").Value = -6693.4375
$ws.Range("H68").Value = 3274.875
$ws.Range("I68").Value = 1717.6
$ws.Range("J68").Value = 3497.3428
$ws.Range("K68").Value = 5152.799999999999
$ws.Range("L68").Value = 10492.0284
$ws.Range("M68").Value = -4341.799999999999
$ws.Range("N68").Value = -12114.0284
$ws.Range("H69").Value = 4799.95
$ws.Range("I69").Value = 3000
$ws.Range("K69").Value = 9000
$ws.Range("M69").Value = -8189
$ws.Range("H71").Value = 3274.875
$ws.Range("I71").Value = 1717.6
$ws.Range("J71").Value = 3497.3428
$ws.Range("K71").Value = 15458.4
$ws.Range("L71").Value = 31476.0852
$ws.Range("M71").Value = -11402.4
$ws.Range("N71").Value = -39588.0852
$ws.Range("H72").Value = 4799.95
$ws.Range("I72").Value = 3000
$ws.Range("K72").Value = 27000
$ws.Range("M72").Value = -22944
$ws.Range("H86").Value = 855.125
$ws.Range("I86").Value = 51.5
$ws.Range("J86").Value = 1123
$ws.Range("K86").Value = 154.5
$ws.Range("L86").Value = 3369
$ws.Range("M86").Value = 1031.5
$ws.Range("N86").Value = -5741
$ws.Range("H89").Value = 855.125
$ws.Range("I89").Value = 51.5
$ws.Range("J89").Value = 1123
$ws.Range("K89").Value = 463.5
$ws.Range("L89").Value = 10107
$ws.Range("M89").Value = 5464.5
$ws.Range("N89").Value = -21963
$ws.Range("H131").Value = 178691.58
$ws.Range("I131").Value = 284514.4
$ws.Range("K131").Value = 853543.2000000001
$ws.Range("M131").Value = -848503.2000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4749.5713
$ws.Range("I70").Value = 5226.25
$ws.Range("J70").Value = 4114
$ws.Range("K70").Value = 5226.25
$ws.Range("L70").Value = 4114
$ws.Range("M70").Value = -4956.25
$ws.Range("N70").Value = -4654
$ws.Range("H73").Value = 4749.5713
$ws.Range("I73").Value = 5226.25
$ws.Range("J73").Value = 4114
$ws.Range("K73").Value = 5226.25
$ws.Range("L73").Value = 4114
$ws.Range("M73").Value = -4290.25
$ws.Range("N73").Value = -5986
$ws.Range("H102").Value = 18820.115
$ws.Range("I102").Value = 23815.852
$ws.Range("K102").Value = 23815.852
$ws.Range("M102").Value = -22193.852
$ws.Range("H126").Value = 4532.75
$ws.Range("I126").Value = 1772.4
$ws.Range("J126").Value = 6504.4287
$ws.Range("K126").Value = 5317.200000000001
$ws.Range("L126").Value = 19513.2861
$ws.Range("M126").Value = -2847.200000000001
$ws.Range("N126").Value = -24453.2861
$ws.Range("H132").Value = 3347.6
$ws.Range("I132").Value = 3423.842
$ws.Range("K132").Value = 10271.526
$ws.Range("M132").Value = -7741.526

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1111.8182
$ws.Range("I22").Value = 680
$ws.Range("J22").Value = 1630
$ws.Range("K22").Value = 680
$ws.Range("L22").Value = 1630
$ws.Range("M22").Value = -385
$ws.Range("N22").Value = -2220
$ws.Range("H27").Value = 1111.8182
$ws.Range("I27").Value = 680
$ws.Range("J27").Value = 1630
$ws.Range("K27").Value = 680
$ws.Range("L27").Value = 1630
$ws.Range("M27").Value = -573
$ws.Range("N27").Value = -1844
$ws.Range("H132").Value = 3246.0334
$ws.Range("I132").Value = 3429.3635
$ws.Range("J132").Value = 2741.875
$ws.Range("K132").Value = 10288.0905
$ws.Range("L132").Value = 8225.625
$ws.Range("M132").Value = -7758.0905
$ws.Range("N132").Value = -13285.625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1358.5454
$ws.Range("I81").Value = 1294.4
$ws.Range("K81").Value = 2588.8
$ws.Range("M81").Value = -1527.8
$ws.Range("H84").Value = 1358.5454
$ws.Range("I84").Value = 1294.4
$ws.Range("K84").Value = 12944
$ws.Range("M84").Value = -7640
